$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Aidan McCarron MCIOB
$ws.Range("A2").Value = "Aidan Mc"

# Row 3 - John Higgins MCIOB
$ws.Range("A3").Value = "John Higgins"
$ws.Range("C3").Value = "Higgins"

# Row 4 - Michael Yohanis MCIOB
$ws.Range("A4").Value = "Michael Yohanis"
$ws.Range("C4").Value = "Yohanis"

# Row 5 - Declan McLogan CMIOSH LL.M
$ws.Range("A5").Value = "Declan Mc"
$ws.Range("C5").Value = "Mc"

# Row 10 - Lee Robert Gray GradIOSH
$ws.Range("A10").Value = "Lee Robert Gray Grad"
$ws.Range("B10").Value = "Lee Robert"
$ws.Range("C10").Value = "Gray Grad"

# Row 11 - Eamonn Laverty. MCIOB
$ws.Range("A11").Value = "Eamonn Laverty"
$ws.Range("C11").Value = "Laverty"

# Row 15 - Sinéad Gorman (she/her)
$ws.Range("A15").Value = "Sinead Gorman"
$ws.Range("B15").Value = "Sinead"
$ws.Range("C15").Value = "Gorman"

# Row 19 - Nina Salandy    BSc. (Hons.) GradIOSH
$ws.Range("A19").Value = "Nina Salandy"
$ws.Range("C19").Value = "Salandy"
